# Swap the order of names in the "Recorded By" column (G) wherever the
# cell reads "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# Cells that contain only "dnasr281@gmail.com" (no "System") are left
# untouched, as are all other columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Use Find/FindNext to locate only the matching cells instead of scanning
# every cell in the sheet (which would needlessly touch/realize blank
# cells in between).
$addresses = New-Object System.Collections.ArrayList

$first = $ws.Cells.Find($oldValue)
if ($first -ne $null) {
    $firstAddr = $first.Address()
    $addresses.Add($firstAddr) | Out-Null

    $current = $ws.Cells.FindNext($first)
    while ($current.Address() -ne $firstAddr) {
        $addresses.Add($current.Address()) | Out-Null
        $current = $ws.Cells.FindNext($current)
    }
}

foreach ($addr in $addresses) {
    $ws.Range($addr).Value2 = $newValue
}
